$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "Creditos_Requisitos"
$ws.Range("E60").Value = "168452/168209"
$ws.Range("E52").Value = "168416/168238"

# Prerequisite lists are now stored as comma/slash separated text, so the
# columns that hold them (and the rest of the used range) are switched to
# the Text number format to keep them from being re-interpreted.
$ws.Range("A1:G63").NumberFormat = "@"
$ws.Range("M17").NumberFormat = "@"

# Column E ("Prerrequisitos") was manually widened (and no longer auto-fit).
$ws.Columns("E").ColumnWidth = 18.75

# View state: zoom + selection as left by the editing session.
$excel.ActiveWindow.Zoom = 132
$null = $ws.Range("C24").Select()
